# Add a new "2022-Q4" sheet (with fund holding data) right after "总计"
# and before "2022-Q3", and add a corresponding summary row to "总计".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert the new worksheet "2022-Q4" right before "2022-Q3"
# ---------------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$newSheet = $wb.Worksheets.Add($q3)
$newSheet.Name = "2022-Q4"

# A pristine, never-formatted cell (used below as a format donor so that
# text values we stamp with a leading apostrophe don't keep a stray
# "quote prefix" style once pasted into the new sheet).
$blank = $q3.Range("A1")

# ---------------------------------------------------------------------------
# 2. Populate the header row
# ---------------------------------------------------------------------------
$headers = @("基金代码","基金名称","基金规模","股票总仓位","仓位占比","持有市值(亿元)","仓位排名")
for ($c = 0; $c -lt $headers.Length; $c++) {
    $newSheet.Cells.Item(1, $c + 2).Value = $headers[$c]
}
$newSheet.Range("B1:H1").Style = $q3.Range("B1").Style

# ---------------------------------------------------------------------------
# 3. Populate the fund rows (basic code|name|size|position|ratio|value|rank)
# ---------------------------------------------------------------------------
$data = @"
012578|富国红利混合A|13.52|93.55|3.58|0.4840|3
400032|东方主题精选混合|12.76|92.86|3.54|0.4517|6
010330|东吴兴享成长混合A|7.86|82.94|5.50|0.4323|3
400003|东方精选混合|9.59|92.83|3.59|0.3443|10
001702|东方创新科技混合|6.83|93.04|4.27|0.2916|7
016097|东吴兴弘一年持有期混合A|4.92|66.01|5.09|0.2504|1
012850|中融低碳经济3个月持有期混合A|5.54|88.38|4.09|0.2266|5
010751|宝盈优质成长混合A|4.33|93.88|4.90|0.2122|5
005296|南华丰淳混合A|4.01|89.14|4.81|0.1929|4
161601|融通新蓝筹混合|11.21|70.40|1.20|0.1345|8
001543|宝盈新锐灵活配置混合A|2.28|91.59|4.82|0.1099|9
011160|富国质量成长6个月持有期混合A|3.97|89.73|2.75|0.1092|9
012579|富国红利混合C|3.04|93.55|3.58|0.1088|3
400001|东方龙混合|2.51|88.26|3.65|0.0916|10
011462|东吴兴享成长混合C|1.17|82.94|5.50|0.0644|3
011458|东方鑫享价值成长一年持有期混合A|2.18|87.88|2.86|0.0623|9
014352|东方创新成长混合A|1.27|89.61|3.96|0.0503|7
012851|中融低碳经济3个月持有期混合C|1.19|88.38|4.09|0.0487|5
016098|东吴兴弘一年持有期混合C|0.76|66.01|5.09|0.0387|1
011459|东方鑫享价值成长一年持有期混合C|1.29|87.88|2.86|0.0369|9
010752|宝盈优质成长混合C|0.70|93.88|4.90|0.0343|5
005297|南华丰淳混合C|0.33|89.14|4.81|0.0159|4
002955|融通新趋势灵活配置混合|0.74|90.61|1.61|0.0119|8
007578|宝盈新锐灵活配置混合C|0.21|91.59|4.82|0.0101|9
015382|东方兴瑞趋势领航混合C|0.27|85.64|3.41|0.0092|7
014353|东方创新成长混合C|0.17|89.61|3.96|0.0067|7
015381|东方兴瑞趋势领航混合A|0.19|85.64|3.41|0.0065|7
015575|宝盈新能源产业混合C|0.13|92.50|5.03|0.0065|7
015574|宝盈新能源产业混合A|0.12|92.50|5.03|0.0060|7
010646|融通价值趋势混合A|0.41|74.27|1.26|0.0052|10
003670|中融物联网主题灵活配置混合|0.13|92.35|3.73|0.0048|4
011161|富国质量成长6个月持有期混合C|0.15|89.73|2.75|0.0041|9
010647|融通价值趋势混合C|0.10|74.27|1.26|0.0013|10
"@

$lines = $data -split "`n"
$rowIndex = 2
foreach ($line in $lines) {
    $fields = $line -split "\|"
    $code = $fields[0]
    $fundName = $fields[1]
    $size = $fields[2]
    $position = $fields[3]
    $ratio = $fields[4]
    $value = $fields[5]
    $rank = [int]$fields[6]

    $newSheet.Cells.Item($rowIndex, 1).Value = $rowIndex - 2
    $newSheet.Cells.Item($rowIndex, 2).Value = "'" + $code
    $newSheet.Cells.Item($rowIndex, 3).Value = $fundName
    $newSheet.Cells.Item($rowIndex, 4).Value = "'" + $size
    $newSheet.Cells.Item($rowIndex, 5).Value = "'" + $position
    $newSheet.Cells.Item($rowIndex, 6).Value = "'" + $ratio
    $newSheet.Cells.Item($rowIndex, 7).Value = "'" + $value
    $newSheet.Cells.Item($rowIndex, 8).Value = $rank

    $rowIndex = $rowIndex + 1
}
$lastRow = $rowIndex - 1

# Column A (row index) uses the same bordered/centered style as the other
# report sheets.
$newSheet.Range("A2:A" + $lastRow).Style = $q3.Range("A2").Style

# The text fields above were entered with a leading apostrophe so Excel
# keeps them as text (e.g. "012578" rather than numeric 12578). That also
# stamps the cells with a "quote prefix" style; strip that back out by
# re-pasting the (default) formatting from a pristine, untouched cell so
# the cells end up with plain/no explicit style, matching the rest of the
# workbook.
$blank.Copy()
$newSheet.Range("B2:B" + $lastRow).PasteSpecial(-4122)
$newSheet.Range("D2:G" + $lastRow).PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 4. Add the corresponding summary row to "总计"
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()
$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q4"
$total.Cells.Item(2, 3).Value = 33
$total.Cells.Item(2, 4).Value = 3.86

# Give the new A2 the same style as the rest of column A, and make sure
# B2:D2 end up with the plain/no-style formatting used by the other rows.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
